$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: add Mitigated/Resolved date, change ROAM to "R", update How text
$ws.Range("F3").Value = "4/2/2023"
$ws.Range("G3").Value = "R"
$ws.Range("H3").Value = "I didn't have time to work but I will still be able to finish on time"

# Row 4: new risk entry
$ws.Range("A4").Value = "Troubles with physics model will delay other parts"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 13
$ws.Range("E4").Value = "4/5/2023"
$ws.Range("G4").Value = "O"
$ws.Range("H4").Value = "Need a plan to mitigate this risk"

# Widen column H and move the active selection to H4
$ws.Columns.Item(8).ColumnWidth = 50
$ws.Range("H4").Select()

# Let the P/I scatter chart pick up the newly added data point
$excel.Calculate()
$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh()
